# Insert a new data row before row 426 (this pushes the existing rows
# 426..489 down to 427..490), then populate the newly inserted row 426
# with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 426..489 down by one row.
$ws.Range("A426").EntireRow.Insert()

# Fill in the new record on row 426.
$ws.Range("A426").Value = 4
$ws.Range("B426").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C426").Value = "Los Lagos"
$ws.Range("D426").Value = 44984
$ws.Range("E426").Value = 10
$ws.Range("F426").Value = 100112008
$ws.Range("G426").Value = "Coliflor"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 250
$ws.Range("K426").Value = 1700
$ws.Range("L426").Value = 1700
$ws.Range("M426").Value = 1700
$ws.Range("N426").Value = "$/unidad"
$ws.Range("O426").Value = "Región Metropolitana"
$ws.Range("P426").Value = 1700
$ws.Range("Q426").Value = 1
$ws.Range("R426").Value = "Hortaliza"
